# Add 2022-Q1 data: insert a new "2022-Q1" sheet (with one fund holding row)
# between "2021-Q4" and "总计", and add a corresponding summary row to "总计".

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: Remove the existing "总计" sheet. We will recreate it further below,
# after the new "2022-Q1" sheet has been created, so that sheetId allocation
# ends up with 2022-Q1 = 6 and 总计 = 7 (matching creation order).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Delete()

# ---------------------------------------------------------------------------
# Step 2: Create the new "2022-Q1" sheet right after "2021-Q4".
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q1Sheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4Sheet)
$q1Sheet.Name = "2022-Q1"

$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"
$q1Sheet.Range("B1:H1").Style = $q4Sheet.Range("B1:H1").Style

$q1Sheet.Range("A2").Value = 0
$q1Sheet.Range("A2").Style = $q4Sheet.Range("A2").Style
# Fund code keeps its leading zero, so it must be forced to stay text.
$q1Sheet.Range("B2").Value = "'001118"
$q1Sheet.Range("C2").Value = "华宝事件驱动混合"
$q1Sheet.Range("D2").Value = 6.29
$q1Sheet.Range("E2").Value = 92.73
$q1Sheet.Range("F2").Value = 3.06
$q1Sheet.Range("G2").Value = 0.1925
$q1Sheet.Range("H2").Value = 10

# ---------------------------------------------------------------------------
# Step 3: Recreate the "总计" sheet after "2022-Q1" and rewrite its table,
# with the new "2022-Q1" row on top and all other quarters shifted down.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1Sheet)
$totalSheet.Name = "总计"

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$rows = @(
    @{ idx = 0; date = "2022-Q1"; count = 1;  value = 0.19 },
    @{ idx = 1; date = "2021-Q4"; count = 2;  value = 0.48 },
    @{ idx = 2; date = "2021-Q3"; count = 2;  value = 0.24 },
    @{ idx = 3; date = "2021-Q2"; count = 5;  value = 1.54 },
    @{ idx = 4; date = "2021-Q1"; count = 13; value = 1.66 },
    @{ idx = 5; date = "2020-Q4"; count = 29; value = 20.63 }
)

$r = 2
foreach ($row in $rows) {
    $totalSheet.Range("A$r").Value = $row.idx
    $totalSheet.Range("B$r").Value = $row.date
    $totalSheet.Range("C$r").Value = $row.count
    $totalSheet.Range("D$r").Value = $row.value
    $r = $r + 1
}
